$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 96 and 97 (columns F:V hold the match data; A:E are the
#     shared index/metadata columns and stay put) ---
$row96 = $ws.Range("F96:V96").Value2
$row97 = $ws.Range("F97:V97").Value2
$ws.Range("F96:V96").Value2 = $row97
$ws.Range("F97:V97").Value2 = $row96

# --- Swap rows 222 and 223 the same way ---
$row222 = $ws.Range("F222:V222").Value2
$row223 = $ws.Range("F223:V223").Value2
$ws.Range("F222:V222").Value2 = $row223
$ws.Range("F223:V223").Value2 = $row222

# --- Append new row 224 (Coquimbo vs Union La Calera) ---
# Seed the new row from row 223's formats so the A/E column styles (bold
# bordered index column, datetime-formatted match-date column) match the
# rest of the table.
$ws.Range("A223:V223").Copy()
$ws.Range("A224:V224").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B224/C224/D224 are identical to row 223's values ("chile",
# "primera-division", "2023"); paste them as values so the numeric-looking
# "2023" stays a text cell instead of being re-interpreted as a number.
$ws.Range("B223:D223").Copy()
$ws.Range("B224:D224").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("A224").Value2 = 223
$ws.Range("E224").Value2 = 45262
$ws.Range("F224").Value2 = "Coquimbo"
$ws.Range("G224").Value2 = 1
$ws.Range("H224").Value2 = "Union La Calera"
$ws.Range("I224").Value2 = 2
$ws.Range("J224").Value2 = 2.23
$ws.Range("K224").Value2 = "25/11/2023 00:42"
$ws.Range("L224").Value2 = 2.27
$ws.Range("M224").Value2 = "01/12/2023 23:51"
$ws.Range("N224").Value2 = 3.4
$ws.Range("O224").Value2 = "25/11/2023 00:42"
$ws.Range("P224").Value2 = 3.4
$ws.Range("Q224").Value2 = "01/12/2023 23:51"
$ws.Range("R224").Value2 = 3.33
$ws.Range("S224").Value2 = "25/11/2023 00:42"
$ws.Range("T224").Value2 = 3.34
$ws.Range("U224").Value2 = "01/12/2023 23:54"
$ws.Range("V224").Value2 = "https://www.betexplorer.com/football/chile/primera-division/coquimbo-union-la-calera/d6W7rPb4/"

Write-Output "edit applied"
